$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "36.643.41"
Set-TextValue $ws "E2" "  +3.14%  "
Set-TextValue $ws "D3" "2.056.32"
Set-TextValue $ws "E3" "  +8.80%  "
Set-TextValue $ws "D4" "0.998"
Set-TextValue $ws "E4" "  -0.30%  "
Set-TextValue $ws "D5" "247.30"
Set-TextValue $ws "E5" "  +0.54%  "
Set-TextValue $ws "D6" "0.663"
Set-TextValue $ws "E6" "  -3.83%  "
Set-TextValue $ws "D7" "0.998"
Set-TextValue $ws "E7" "  -0.17%  "
Set-TextValue $ws "D8" "45.13"
Set-TextValue $ws "E8" "  +4.76%  "
Set-TextValue $ws "D9" "60.33"
Set-TextValue $ws "E9" "  +6.66%  "
Set-TextValue $ws "D10" "0.363"
Set-TextValue $ws "E10" "  +2.16%  "
Set-TextValue $ws "D11" "0.0721"
Set-TextValue $ws "E11" "  -4.07%  "
Set-TextValue $ws "D12" "0.0986"
Set-TextValue $ws "E12" "  +0.30%  "
Set-TextValue $ws "D13" "14.55"
Set-TextValue $ws "E13" "  -0.11%  "
Set-TextValue $ws "D14" "2.334.93"
Set-TextValue $ws "E14" "  +7.39%  "
Set-TextValue $ws "D15" "0.810"
Set-TextValue $ws "E15" "  +2.35%  "
Set-TextValue $ws "D16" "2.031.35"
Set-TextValue $ws "E16" "  +6.94%  "
Set-TextValue $ws "D17" "4.91"
Set-TextValue $ws "E17" "  -2.22%  "
Set-TextValue $ws "D18" "36.571.85"
Set-TextValue $ws "E18" "  +2.92%  "
Set-TextValue $ws "D19" "71.38"
Set-TextValue $ws "E19" "  -2.76%  "
Set-TextValue $ws "D20" "0.0₃0815"
Set-TextValue $ws "E20" "  -1.59%  "
Set-TextValue $ws "D21" "237.35"
Set-TextValue $ws "E21" "  -3.59%  "
Set-TextValue $ws "D22" "12.62"
Set-TextValue $ws "E22" "  -3.06%  "
Set-TextValue $ws "D23" "4.92"
Set-TextValue $ws "E23" "  -4.77%  "
Set-TextValue $ws "D25" "2.47"
Set-TextValue $ws "E25" "  -7.24%  "
Set-TextValue $ws "D26" "169.23"
Set-TextValue $ws "E26" "  +1.88%  "
Set-TextValue $ws "B27" "EthereumClassic"
Set-TextValue $ws "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D27" "20.12"
Set-TextValue $ws "E27" "  +9.74%  "
Set-TextValue $ws "B28" "Cosmos"
Set-TextValue $ws "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D28" "8.79"
Set-TextValue $ws "E28" "  +1.75%  "
Set-TextValue $ws "D29" "1.95"
Set-TextValue $ws "E29" "  -8.65%  "
Set-TextValue $ws "E30" "  -4.80%  "
Set-TextValue $ws "D31" "21.72"
Set-TextValue $ws "E31" "  +51.05%  "
Set-TextValue $ws "D32" "4.37"
Set-TextValue $ws "E32" "  -0.91%  "
Set-TextValue $ws "D33" "0.0581"
Set-TextValue $ws "E33" "  -4.34%  "
Set-TextValue $ws "B34" "Kaspa"
Set-TextValue $ws "C34" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D34" "0.0889"
Set-TextValue $ws "E34" "  +18.78%  "
Set-TextValue $ws "B35" "BinanceUSD"
Set-TextValue $ws "C35" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  -0.22%  "
Set-TextValue $ws "E36" "  +1.12%  "
Set-TextValue $ws "D37" "2.28"
Set-TextValue $ws "E37" "  +17.46%  "
Set-TextValue $ws "D38" "3.99"
Set-TextValue $ws "E38" "  -6.36%  "
Set-TextValue $ws "D39" "0.867"
Set-TextValue $ws "E39" "  +1.75%  "
Set-TextValue $ws "D40" "1.32"
Set-TextValue $ws "E40" "  -10.64%  "
Set-TextValue $ws "D41" "0.0215"
Set-TextValue $ws "E41" "  -6.01%  "
Set-TextValue $ws "D42" "96.43"
Set-TextValue $ws "E42" "  -2.62%  "
Set-TextValue $ws "D43" "1.12"
Set-TextValue $ws "E43" "  +3.18%  "
Set-TextValue $ws "D44" "2.78"
Set-TextValue $ws "E44" "  +15.66%  "
Set-TextValue $ws "D45" "16.00"
Set-TextValue $ws "E45" "  -5.64%  "
Set-TextValue $ws "D46" "1.317.04"
Set-TextValue $ws "E46" "  +0.31%  "
Set-TextValue $ws "D47" "0.0814"
Set-TextValue $ws "E47" "  +0.28%  "
Set-TextValue $ws "D48" "2.81"
Set-TextValue $ws "E48" "  +2.70%  "
Set-TextValue $ws "D49" "2.226.90"
Set-TextValue $ws "E49" "  +7.33%  "
Set-TextValue $ws "E50" "  -5.90%  "
Set-TextValue $ws "D51" "3.84"
Set-TextValue $ws "E51" "  +15.80%  "
